$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has (rows 2-5):
#   2: Paper (Ori)
#   3: Paper (rerun)
#   4: MiniInceptionTime (keras)     98.97 98.97 99.45 98.96 6261
#   5: MiniInceptionTime (pytorch)   x     x     x     x     x
#
# New layout after this edit (rows 2-6):
#   2: Paper (Ori)                          [unchanged]
#   3: Paper (rerun – pytorch)              [unchanged values, relabeled]
#   4: Paper (rerun – keras)                98.85 98.85 99.43 98.83 11065  [new row]
#   5: MiniInceptionTime (keras – 2560)     98.97 98.97 99.45 98.96 6261   [old row4 values, relabeled]
#   6: MiniInceptionTime (keras – 1280)     98.85 98.85 99.39 98.83 6261   [relabeled + new values]

# Insert a new row above the old row 4, shifting the two "MiniInceptionTime" rows down.
$ws.Rows("4").Insert()

# Row 3: relabel "Paper (rerun)" -> "Paper (rerun – pytorch)" (values stay the same).
$ws.Range("A3").Value = "Paper (rerun – pytorch)"

# Row 4 (brand new row): "Paper (rerun – keras)" with its own results.
$ws.Range("A4").Value = "Paper (rerun – keras)"
$ws.Range("B4").Value = 98.85
$ws.Range("C4").Value = 98.85
$ws.Range("D4").Value = 99.43
$ws.Range("E4").Value = 98.83
$ws.Range("F4").Value = 11065

# Row 5 (previously row 4): relabel to the keras-2560 variant; numbers are unchanged.
$ws.Range("A5").Value = "MiniInceptionTime (keras – 2560)"

# Row 6 (previously row 5): relabel to the keras-1280 variant and fill in real numbers
# (replacing the old placeholder "x" values).
$ws.Range("A6").Value = "MiniInceptionTime (keras – 1280)"
$ws.Range("B6").Value = 98.85
$ws.Range("C6").Value = 98.85
$ws.Range("D6").Value = 99.39
$ws.Range("E6").Value = 98.83
$ws.Range("F6").Value = 6261

# Column A needs to be a bit wider to fit the new, longer labels.
$ws.Columns("A").ColumnWidth = 30

# Restore the (unrelated) leftover cell selection recorded in the sheet view.
$ws.Range("J7").Select()
